$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.155.34"
$ws.Range("E2").Value = "  -2.17%  "
$ws.Range("D3").Value = "1.571.68"
$ws.Range("E3").Value = "  -1.46%  "
$ws.Range("E4").Value = "  -0.33%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.12%  "
$ws.Range("E6").Value = "  -2.60%  "
$ws.Range("E7").Value = "  -0.31%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.33"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.17%  "
$ws.Range("E9").Value = "  -2.18%  "
$ws.Range("E10").Value = "  -0.47%  "
$ws.Range("E11").Value = "  -0.48%  "
$ws.Range("D12").Value = "1.792.26"
$ws.Range("E12").Value = "  -1.61%  "
$ws.Range("D13").Value = "1.568.51"
$ws.Range("E13").Value = "  -1.71%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.77"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.14%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.519"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.52%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "62.70"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.18%  "
$ws.Range("D17").Value = "27.166.52"
$ws.Range("E17").Value = "  -2.13%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "214.76"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.20%  "
$ws.Range("E19").Value = "  -1.18%  "
$ws.Range("E20").Value = "  -1.45%  "
$ws.Range("E21").Value = "  -0.27%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.14"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.57%  "
$ws.Range("E23").Value = "  -3.60%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.14%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.47"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.66"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -6.44%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.94"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.50%  "
$ws.Range("B28").Value = "Stellar"
$ws.Range("C28").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.104"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.36%  "
$ws.Range("B29").Value = "BinanceUSD"
$ws.Range("C29").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.29%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.14"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.70%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0464"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.08%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.18"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.83%  "
$ws.Range("D33").Value = "1.399.31"
$ws.Range("E33").Value = "  +1.32%  "
$ws.Range("E34").Value = "  -1.84%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.56"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.99%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.946"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.95%  "
$ws.Range("E37").Value = "  -2.30%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0166"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.87%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.816"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.56%  "
$ws.Range("E40").Value = "  -3.59%  "
$ws.Range("E41").Value = "  -0.28%  "
$ws.Range("E42").Value = "  +1.63%  "
$ws.Range("E43").Value = "  +2.86%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.35"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.18%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "63.70"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.26%  "
$ws.Range("B46").Value = "MXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.17"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.15%  "
$ws.Range("D47").Value = "1.706.11"
$ws.Range("E47").Value = "  -1.53%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "85.91"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.41%  "
$ws.Range("D49").Value = "0.0₇0981"
$ws.Range("E49").Value = "  -3.02%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0951"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.67%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0493"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.50%  "
